$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name "Feuille1" -> "Sheet")
$ws.Name = "Sheet"

# Remove the now-optional "numero_article" header (D1)
$ws.Range("D1").ClearContents()

$wb.Save()
